# Update raw market-data values pulled by the scheduled Sheets runner.
# Each sheet stores currentAveragePrice/.../LeveProfit* columns (H:N) as plain
# numbers (no formulas) pulled from an external price feed, so this just
# overwrites the specific cells whose upstream source values changed.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(11, 8).Value = 60.333332  # H11: 60.416668 -> 60.333332
$ws.Cells.Item(11, 9).Value = 60.333332  # I11: 60.416668 -> 60.333332
$ws.Cells.Item(11, 11).Value = 60.333332  # K11: 60.416668 -> 60.333332
$ws.Cells.Item(11, 13).Value = 79.666668  # M11: 79.583332 -> 79.666668
$ws.Cells.Item(33, 8).Value = 210.5  # H33: 222.6 -> 210.5
$ws.Cells.Item(33, 9).Value = 228.9  # I33: 248.625 -> 228.9
$ws.Cells.Item(33, 11).Value = 228.9  # K33: 248.625 -> 228.9
$ws.Cells.Item(33, 13).Value = 0.09999999999999432  # M33: -19.625 -> 0.09999999999999432
$ws.Cells.Item(106, 8).Value = 37039572  # H106: 41669490 -> 37039572
$ws.Cells.Item(106, 9).Value = 47619450  # I106: 55555984 -> 47619450
$ws.Cells.Item(106, 11).Value = 47619450  # K106: 55555984 -> 47619450
$ws.Cells.Item(106, 13).Value = -47618819  # M106: -55555353 -> -47618819
$ws.Cells.Item(137, 8).Value = 4428.2856  # H137: 4799.6 -> 4428.2856
$ws.Cells.Item(137, 9).Value = 3199.6  # I137: 3499.5 -> 3199.6
$ws.Cells.Item(137, 10).Value = 7500  # J137: 10000 -> 7500
$ws.Cells.Item(137, 11).Value = 9598.799999999999  # K137: 10498.5 -> 9598.799999999999
$ws.Cells.Item(137, 12).Value = 22500  # L137: 30000 -> 22500
$ws.Cells.Item(137, 13).Value = -7048.799999999999  # M137: -7948.5 -> -7048.799999999999
$ws.Cells.Item(137, 14).Value = -27600  # N137: -35100 -> -27600
$ws.Cells.Item(138, 8).Value = 8408.637000000001  # H138: 8156.543 -> 8408.637000000001
$ws.Cells.Item(138, 10).Value = 8644.678  # J138: 8363 -> 8644.678
$ws.Cells.Item(138, 12).Value = 25934.034  # L138: 25089 -> 25934.034
$ws.Cells.Item(138, 14).Value = -36214.034  # N138: -35369 -> -36214.034

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(19, 8).Value = 1999.5  # H19: 2000 -> 1999.5
$ws.Cells.Item(19, 9).Value = 1999.5  # I19: 2000 -> 1999.5
$ws.Cells.Item(19, 11).Value = 1999.5  # K19: 2000 -> 1999.5
$ws.Cells.Item(19, 13).Value = -1770.5  # M19: -1771 -> -1770.5
$ws.Cells.Item(32, 8).Value = 1630.3594  # H32: 1697.5469 -> 1630.3594
$ws.Cells.Item(32, 9).Value = 1105.5667  # I32: 1115.8306 -> 1105.5667
$ws.Cells.Item(32, 10).Value = 9502.25  # J32: 8561.799999999999 -> 9502.25
$ws.Cells.Item(32, 11).Value = 1105.5667  # K32: 1115.8306 -> 1105.5667
$ws.Cells.Item(32, 12).Value = 9502.25  # L32: 8561.799999999999 -> 9502.25
$ws.Cells.Item(32, 13).Value = -818.5667000000001  # M32: -828.8306 -> -818.5667000000001
$ws.Cells.Item(32, 14).Value = -10076.25  # N32: -9135.799999999999 -> -10076.25
$ws.Cells.Item(43, 8).Value = 8955.200000000001  # H43: 8944 -> 8955.200000000001
$ws.Cells.Item(43, 10).Value = 9569  # J43: 9758.666999999999 -> 9569
$ws.Cells.Item(43, 12).Value = 9569  # L43: 9758.666999999999 -> 9569
$ws.Cells.Item(43, 14).Value = -10195  # N43: -10384.667 -> -10195
$ws.Cells.Item(97, 8).Value = 1093.2858  # H97: 1050.409 -> 1093.2858
$ws.Cells.Item(97, 9).Value = 924.6842  # I97: 885.95 -> 924.6842
$ws.Cells.Item(97, 11).Value = 924.6842  # K97: 885.95 -> 924.6842
$ws.Cells.Item(97, 13).Value = -428.6842  # M97: -389.95 -> -428.6842
$ws.Cells.Item(122, 8).Value = 328805.72  # H122: 349293.56 -> 328805.72
$ws.Cells.Item(122, 9).Value = 613385.75  # I122: 689934 -> 613385.75
$ws.Cells.Item(122, 11).Value = 1840157.25  # K122: 2069802 -> 1840157.25
$ws.Cells.Item(122, 13).Value = -1837707.25  # M122: -2067352 -> -1837707.25
$ws.Cells.Item(132, 8).Value = 4705.8115  # H132: 4773.231 -> 4705.8115
$ws.Cells.Item(132, 10).Value = 5112  # J132: 5256.8887 -> 5112
$ws.Cells.Item(132, 12).Value = 15336  # L132: 15770.6661 -> 15336
$ws.Cells.Item(132, 14).Value = -20396  # N132: -20830.6661 -> -20396

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(5, 8).Value = 9600  # H5: 15000 -> 9600
$ws.Cells.Item(5, 9).Value = 9600  # I5: 15000 -> 9600
$ws.Cells.Item(5, 11).Value = 9600  # K5: 15000 -> 9600
$ws.Cells.Item(5, 13).Value = -9487  # M5: -14887 -> -9487
$ws.Cells.Item(20, 8).Value = 1811.381  # H20: 1777.8636 -> 1811.381
$ws.Cells.Item(20, 9).Value = 1176  # I20: 1169.625 -> 1176
$ws.Cells.Item(20, 11).Value = 1176  # K20: 1169.625 -> 1176
$ws.Cells.Item(20, 13).Value = -929  # M20: -922.625 -> -929
$ws.Cells.Item(134, 8).Value = 3009.5322  # H134: 3053.7212 -> 3009.5322
$ws.Cells.Item(134, 9).Value = 2478.0667  # I134: 2558.1162 -> 2478.0667
$ws.Cells.Item(134, 10).Value = 4416.353  # J134: 4237.6665 -> 4416.353
$ws.Cells.Item(134, 11).Value = 7434.2001  # K134: 7674.348599999999 -> 7434.2001
$ws.Cells.Item(134, 12).Value = 13249.059  # L134: 12712.9995 -> 13249.059
$ws.Cells.Item(134, 13).Value = -4899.2001  # M134: -5139.348599999999 -> -4899.2001
$ws.Cells.Item(134, 14).Value = -18319.059  # N134: -17782.9995 -> -18319.059

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 41671710  # H31: 45459932 -> 41671710
$ws.Cells.Item(31, 9).Value = 250001150  # I31: 500000500 -> 250001150
$ws.Cells.Item(31, 10).Value = 5823.2  # J31: 5873.4 -> 5823.2
$ws.Cells.Item(31, 11).Value = 250001150  # K31: 500000500 -> 250001150
$ws.Cells.Item(31, 12).Value = 5823.2  # L31: 5873.4 -> 5823.2
$ws.Cells.Item(31, 13).Value = -250000855  # M31: -500000205 -> -250000855
$ws.Cells.Item(31, 14).Value = -6413.2  # N31: -6463.4 -> -6413.2
$ws.Cells.Item(34, 8).Value = 41671710  # H34: 45459932 -> 41671710
$ws.Cells.Item(34, 9).Value = 250001150  # I34: 500000500 -> 250001150
$ws.Cells.Item(34, 10).Value = 5823.2  # J34: 5873.4 -> 5823.2
$ws.Cells.Item(34, 11).Value = 250001150  # K34: 500000500 -> 250001150
$ws.Cells.Item(34, 12).Value = 5823.2  # L34: 5873.4 -> 5823.2
$ws.Cells.Item(34, 13).Value = -250000948  # M34: -500000298 -> -250000948
$ws.Cells.Item(34, 14).Value = -6227.2  # N34: -6277.4 -> -6227.2
$ws.Cells.Item(41, 8).Value = 134500  # H41: 153000 -> 134500
$ws.Cells.Item(41, 9).Value = 79333.336  # I41: 79500 -> 79333.336
$ws.Cells.Item(41, 11).Value = 79333.336  # K41: 79500 -> 79333.336
$ws.Cells.Item(41, 13).Value = -78905.336  # M41: -79072 -> -78905.336
$ws.Cells.Item(62, 8).Value = 40749.75  # H62: 46285.43 -> 40749.75
$ws.Cells.Item(62, 10).Value = 45999.715  # J62: 53333 -> 45999.715
$ws.Cells.Item(62, 12).Value = 45999.715  # L62: 53333 -> 45999.715
$ws.Cells.Item(62, 14).Value = -47247.715  # N62: -54581 -> -47247.715
$ws.Cells.Item(65, 8).Value = 40749.75  # H65: 46285.43 -> 40749.75
$ws.Cells.Item(65, 10).Value = 45999.715  # J65: 53333 -> 45999.715
$ws.Cells.Item(65, 12).Value = 229998.575  # L65: 266665 -> 229998.575
$ws.Cells.Item(65, 14).Value = -236238.575  # N65: -272905 -> -236238.575
$ws.Cells.Item(99, 8).Value = 8511.218000000001  # H99: 9193.522999999999 -> 8511.218000000001
$ws.Cells.Item(99, 9).Value = 14850.556  # I99: 16556.875 -> 14850.556
$ws.Cells.Item(99, 10).Value = 4435.9287  # J99: 4662.231 -> 4435.9287
$ws.Cells.Item(99, 11).Value = 14850.556  # K99: 16556.875 -> 14850.556
$ws.Cells.Item(99, 12).Value = 4435.9287  # L99: 4662.231 -> 4435.9287
$ws.Cells.Item(99, 13).Value = -13352.556  # M99: -15058.875 -> -13352.556
$ws.Cells.Item(99, 14).Value = -7431.9287  # N99: -7658.231 -> -7431.9287
$ws.Cells.Item(126, 8).Value = 8511.218000000001  # H126: 9193.522999999999 -> 8511.218000000001
$ws.Cells.Item(126, 9).Value = 14850.556  # I126: 16556.875 -> 14850.556
$ws.Cells.Item(126, 10).Value = 4435.9287  # J126: 4662.231 -> 4435.9287
$ws.Cells.Item(126, 11).Value = 44551.66800000001  # K126: 49670.625 -> 44551.66800000001
$ws.Cells.Item(126, 12).Value = 13307.7861  # L126: 13986.693 -> 13307.7861
$ws.Cells.Item(126, 13).Value = -42081.66800000001  # M126: -47200.625 -> -42081.66800000001
$ws.Cells.Item(126, 14).Value = -18247.7861  # N126: -18926.693 -> -18247.7861
$ws.Cells.Item(132, 8).Value = 63504028  # H132: 57982004 -> 63504028
$ws.Cells.Item(132, 9).Value = 74076750  # I132: 66669156 -> 74076750
$ws.Cells.Item(132, 11).Value = 222230250  # K132: 200007468 -> 222230250
$ws.Cells.Item(132, 13).Value = -222227720  # M132: -200004938 -> -222227720
$ws.Cells.Item(134, 8).Value = 981.13336  # H134: 975.625 -> 981.13336
$ws.Cells.Item(134, 10).Value = 1000  # J134: 974.8333 -> 1000
$ws.Cells.Item(134, 12).Value = 3000  # L134: 2924.4999 -> 3000
$ws.Cells.Item(134, 14).Value = -8070  # N134: -7994.4999 -> -8070

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(32, 8).Value = 82900.2  # H32: 53124.75 -> 82900.2
$ws.Cells.Item(32, 9).Value = 102725.25  # I32: 60199.715 -> 102725.25
$ws.Cells.Item(32, 11).Value = 308175.75  # K32: 180599.145 -> 308175.75
$ws.Cells.Item(32, 13).Value = -307892.75  # M32: -180316.145 -> -307892.75
$ws.Cells.Item(122, 8).Value = 854.2857  # H122: 911.4286 -> 854.2857
$ws.Cells.Item(122, 10).Value = 1139.4  # J122: 1299.4 -> 1139.4
$ws.Cells.Item(122, 12).Value = 10254.6  # L122: 11694.6 -> 10254.6
$ws.Cells.Item(122, 14).Value = -15154.6  # N122: -16594.6 -> -15154.6

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 15878849  # H70: 7942780.5 -> 15878849
$ws.Cells.Item(70, 9).Value = 23813274  # I70: 11909171 -> 23813274
$ws.Cells.Item(70, 10).Value = 9999  # J70: 9999.5 -> 9999
$ws.Cells.Item(70, 11).Value = 23813274  # K70: 11909171 -> 23813274
$ws.Cells.Item(70, 12).Value = 9999  # L70: 9999.5 -> 9999
$ws.Cells.Item(70, 13).Value = -23813004  # M70: -11908901 -> -23813004
$ws.Cells.Item(70, 14).Value = -10539  # N70: -10539.5 -> -10539
$ws.Cells.Item(73, 8).Value = 15878849  # H73: 7942780.5 -> 15878849
$ws.Cells.Item(73, 9).Value = 23813274  # I73: 11909171 -> 23813274
$ws.Cells.Item(73, 10).Value = 9999  # J73: 9999.5 -> 9999
$ws.Cells.Item(73, 11).Value = 23813274  # K73: 11909171 -> 23813274
$ws.Cells.Item(73, 12).Value = 9999  # L73: 9999.5 -> 9999
$ws.Cells.Item(73, 13).Value = -23812338  # M73: -11908235 -> -23812338
$ws.Cells.Item(73, 14).Value = -11871  # N73: -11871.5 -> -11871
$ws.Cells.Item(97, 8).Value = 502.0909  # H97: 512.3 -> 502.0909
$ws.Cells.Item(97, 10).Value = 560  # J97: 600 -> 560
$ws.Cells.Item(97, 12).Value = 560  # L97: 600 -> 560
$ws.Cells.Item(97, 14).Value = -1552  # N97: -1592 -> -1552
$ws.Cells.Item(132, 8).Value = 4659.222  # H132: 5735.273 -> 4659.222
$ws.Cells.Item(132, 9).Value = 4054.0625  # I132: 4898.5557 -> 4054.0625
$ws.Cells.Item(132, 11).Value = 12162.1875  # K132: 14695.6671 -> 12162.1875
$ws.Cells.Item(132, 13).Value = -9632.1875  # M132: -12165.6671 -> -9632.1875

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(16, 8).Value = 70591130  # H16: 60002532 -> 70591130
$ws.Cells.Item(16, 9).Value = 85715380  # I16: 80001090 -> 85715380
$ws.Cells.Item(16, 10).Value = 11301.333  # J16: 6860.8 -> 11301.333
$ws.Cells.Item(16, 11).Value = 85715380  # K16: 80001090 -> 85715380
$ws.Cells.Item(16, 12).Value = 11301.333  # L16: 6860.8 -> 11301.333
$ws.Cells.Item(16, 13).Value = -85715210  # M16: -80000920 -> -85715210
$ws.Cells.Item(16, 14).Value = -11641.333  # N16: -7200.8 -> -11641.333
$ws.Cells.Item(46, 8).Value = 4567.8335  # H46: 4858.4287 -> 4567.8335
$ws.Cells.Item(46, 9).Value = 2874.8  # I46: 4458.3335 -> 2874.8
$ws.Cells.Item(46, 11).Value = 2874.8  # K46: 4458.3335 -> 2874.8
$ws.Cells.Item(46, 13).Value = -2686.8  # M46: -4270.3335 -> -2686.8
$ws.Cells.Item(136, 8).Value = 3411.44  # H136: 3490.8147 -> 3411.44
$ws.Cells.Item(136, 9).Value = 2928.4443  # I136: 2879.3684 -> 2928.4443
$ws.Cells.Item(136, 10).Value = 4653.4287  # J136: 4943 -> 4653.4287
$ws.Cells.Item(136, 11).Value = 8785.332900000001  # K136: 8638.1052 -> 8785.332900000001
$ws.Cells.Item(136, 12).Value = 13960.2861  # L136: 14829 -> 13960.2861
$ws.Cells.Item(136, 13).Value = -6235.332900000001  # M136: -6088.1052 -> -6235.332900000001
$ws.Cells.Item(136, 14).Value = -19060.2861  # N136: -19929 -> -19060.2861

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(136, 8).Value = 8765.064  # H136: 9080.378000000001 -> 8765.064
$ws.Cells.Item(136, 9).Value = 1545.7  # I136: 1785 -> 1545.7
$ws.Cells.Item(136, 11).Value = 4637.1  # K136: 5355 -> 4637.1
$ws.Cells.Item(136, 13).Value = -2087.1  # M136: -2805 -> -2087.1
